$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1639344262295082
$ws.Range("C2").Value = 0.5942622950819673
$ws.Range("J2").Value = 0.04098360655737705
$ws.Range("P2").Value = 0.110655737704918
$ws.Range("S2").Value = 0.09016393442622951
# Row 3
$ws.Range("B3").Value = 0.006896551724137931
$ws.Range("C3").Value = 0.02068965517241379
$ws.Range("J3").Value = 0.02068965517241379
$ws.Range("P3").Value = 0.7724137931034483
$ws.Range("S3").Value = 0.1793103448275862
# Row 4
$ws.Range("J4").Value = 0.02439024390243903
$ws.Range("P4").Value = 0.7804878048780488
$ws.Range("S4").Value = 0.1951219512195122
# Row 5
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5
# Row 6
$ws.Range("B6").Value = 0.05603448275862069
$ws.Range("D6").Value = 0.008620689655172414
$ws.Range("F6").Value = 0.04741379310344827
$ws.Range("J6").Value = 0.2413793103448276
$ws.Range("O6").Value = 0.008620689655172414
$ws.Range("Q6").Value = 0.1681034482758621
$ws.Range("R6").Value = 0.06896551724137931
$ws.Range("S6").Value = 0.4008620689655172
# Row 7
$ws.Range("B7").Value = 0.09583333333333334
$ws.Range("D7").Value = 0.01666666666666667
$ws.Range("F7").Value = 0.07083333333333333
$ws.Range("J7").Value = 0.1375
$ws.Range("O7").Value = 0.02083333333333333
$ws.Range("Q7").Value = 0.2
$ws.Range("R7").Value = 0.08333333333333333
$ws.Range("S7").Value = 0.375
# Row 8
$ws.Range("B8").Value = 0.06440677966101695
$ws.Range("D8").Value = 0.02033898305084746
$ws.Range("E8").Value = 0.005084745762711864
$ws.Range("F8").Value = 0.05932203389830509
$ws.Range("J8").Value = 0.1152542372881356
$ws.Range("O8").Value = 0.02203389830508475
$ws.Range("Q8").Value = 0.2169491525423729
$ws.Range("R8").Value = 0.09491525423728814
$ws.Range("S8").Value = 0.4016949152542373
# Row 9
$ws.Range("B9").Value = 0.06976744186046512
$ws.Range("D9").Value = 0.01162790697674419
$ws.Range("F9").Value = 0.1511627906976744
$ws.Range("J9").Value = 0.09302325581395349
$ws.Range("O9").Value = 0.03488372093023256
$ws.Range("Q9").Value = 0.1802325581395349
$ws.Range("R9").Value = 0.08139534883720931
$ws.Range("S9").Value = 0.3779069767441861
# Row 10
$ws.Range("B10").Value = 0.08922829581993569
$ws.Range("D10").Value = 0.01929260450160772
$ws.Range("E10").Value = 0.0008038585209003215
$ws.Range("F10").Value = 0.06511254019292605
$ws.Range("J10").Value = 0.09485530546623794
$ws.Range("O10").Value = 0.02090032154340836
$ws.Range("Q10").Value = 0.2564308681672026
$ws.Range("R10").Value = 0.08842443729903537
$ws.Range("S10").Value = 0.364951768488746
# Row 11
$ws.Range("G11").Value = 0.155688622754491
$ws.Range("J11").Value = 0.0718562874251497
$ws.Range("K11").Value = 0.1616766467065868
$ws.Range("L11").Value = 0.6017964071856288
$ws.Range("S11").Value = 0.008982035928143712
# Row 12
$ws.Range("G12").Value = 0.7826086956521739
$ws.Range("J12").Value = 0.1497584541062802
$ws.Range("K12").Value = 0.00966183574879227
$ws.Range("L12").Value = 0.02415458937198068
$ws.Range("S12").Value = 0.03381642512077294
# Row 13
$ws.Range("G13").Value = 0.6470588235294118
$ws.Range("J13").Value = 0.2549019607843137
$ws.Range("S13").Value = 0.09803921568627451
# Row 15
$ws.Range("F15").Value = 0.01345291479820628
$ws.Range("H15").Value = 0.2062780269058296
$ws.Range("I15").Value = 0.05829596412556054
$ws.Range("J15").Value = 0.3183856502242152
$ws.Range("K15").Value = 0.06726457399103139
$ws.Range("M15").Value = 0.004484304932735426
$ws.Range("O15").Value = 0.02690582959641256
$ws.Range("S15").Value = 0.304932735426009
# Row 16
$ws.Range("F16").Value = 0.01234567901234568
$ws.Range("H16").Value = 0.191358024691358
$ws.Range("I16").Value = 0.1111111111111111
$ws.Range("J16").Value = 0.3703703703703703
$ws.Range("K16").Value = 0.1111111111111111
$ws.Range("M16").Value = 0.03703703703703703
$ws.Range("O16").Value = 0.04320987654320987
$ws.Range("S16").Value = 0.1234567901234568
# Row 17
$ws.Range("F17").Value = 0.01769911504424779
$ws.Range("H17").Value = 0.2336283185840708
$ws.Range("I17").Value = 0.06902654867256637
$ws.Range("J17").Value = 0.4
$ws.Range("K17").Value = 0.09734513274336283
$ws.Range("M17").Value = 0.02123893805309734
$ws.Range("N17").Value = 0.001769911504424779
$ws.Range("O17").Value = 0.06017699115044248
$ws.Range("S17").Value = 0.09911504424778761
# Row 18
$ws.Range("F18").Value = 0.0187793427230047
$ws.Range("H18").Value = 0.2206572769953052
$ws.Range("I18").Value = 0.08450704225352113
$ws.Range("J18").Value = 0.3661971830985916
$ws.Range("K18").Value = 0.1173708920187793
$ws.Range("M18").Value = 0.0187793427230047
$ws.Range("O18").Value = 0.06572769953051644
$ws.Range("S18").Value = 0.107981220657277
# Row 19
$ws.Range("F19").Value = 0.0129081245254366
$ws.Range("H19").Value = 0.2528473804100228
$ws.Range("I19").Value = 0.06605922551252848
$ws.Range("J19").Value = 0.3454821564160972
$ws.Range("K19").Value = 0.1161731207289294
$ws.Range("M19").Value = 0.02277904328018223
$ws.Range("N19").Value = 0.002277904328018223
$ws.Range("O19").Value = 0.06302201974183751
$ws.Range("S19").Value = 0.1184510250569476
